# Applies the "landing, noti, red frontend" edit: Word's spell-checker
# (proofErr spellStart/spellEnd) ran over the document, which splits a
# number of runs at word boundaries, and two new SQL snippets were
# appended at the end of the document (about the "skills" table).
#
# NOTE: this COM-interop host only binds PowerShell function parameters
# positionally (named parameters like -Foo bar do not bind), and passing
# a parenthesized "(...)" expression directly as a positional argument
# also misparses - so every multi-part string below is built in a $xml
# variable first, then passed as a bare variable argument.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($index, $innerXml) {
    $para = $d.Paragraphs($index)
    $range = $para.Range
    $xml = "<w:p $wNs>$innerXml</w:p>"
    $range.InsertXML($xml) | Out-Null
}

function Add-ParaXml($innerXml) {
    $lastPara = $d.Paragraphs($d.Paragraphs.Count)
    $lastPara.Range.InsertParagraphAfter() | Out-Null
    $newIndex = $d.Paragraphs.Count
    Set-ParaXml $newIndex $innerXml
}

# 1: "Registro Usuarios"
$xml = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Registro</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Usuarios</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
Set-ParaXml 1 $xml

# 3: "Validacion Usuarios" (paragraph mark keeps u=single + lang=es-ES)
$xml = '<w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Validacion</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> Usuarios</w:t></w:r>'
Set-ParaXml 3 $xml

# 23: '  "POSTS"."datepost" DESC'
$xml = '<w:r><w:t xml:space="preserve">  "POSTS"."</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>datepost</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>" DESC</w:t></w:r>'
Set-ParaXml 23 $xml

# 25: "Información de un usuario"
$xml = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Información</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> de un </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>usuario</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
Set-ParaXml 25 $xml

# 28: "    users.email,            "
$xml = '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>users.email</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">,            </w:t></w:r>'
Set-ParaXml 28 $xml

# 29: "    users.banner,            "
$xml = '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>users.banner</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">,            </w:t></w:r>'
Set-ParaXml 29 $xml

# 30: "    users.photo,            " (keeps lastRenderedPageBreak on first run)
$xml = '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">    </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>users.photo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">,            </w:t></w:r>'
Set-ParaXml 30 $xml

# 32: "    users.location,            "
$xml = '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>users.location</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">,            </w:t></w:r>'
Set-ParaXml 32 $xml

# 33: "    categories.nombre as nombreCAteoria "
$xml = '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>categories.nombre</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> as </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>nombreCAteoria</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
Set-ParaXml 33 $xml

# 35: "INNER JOIN categories ON users.categories_id = categories.id "
$xml = '<w:r><w:t xml:space="preserve">INNER JOIN categories ON </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>users.categories_id</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> = categories.id </w:t></w:r>'
Set-ParaXml 35 $xml

# 38: "Crear Publicacion"
$xml = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Crear</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Publicacion</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
Set-ParaXml 38 $xml

# 39: "INSERT INTO posts (users_id, content, media, type, datepost) VALUES (65,'test 1', 'LinkedIn_logo_initials.png', 'image/png', '2023-10-30 01:20:58')"
# (paragraph mark keeps u=single; every run keeps u=single too)
$xml = '<w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>INSERT INTO posts (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>users_id</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">, content, media, type, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>datepost</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>) VALUES (65,''test 1'', ''LinkedIn_logo_initials.png'', ''image/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>png</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>'', ''2023-10-30 01:20:58'')</w:t></w:r>'
Set-ParaXml 39 $xml

# 41: "Crear comentario"
$xml = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Crear</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>comentario</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
Set-ParaXml 41 $xml

# 42: INSERT INTO COMMENTS (...)
$xml = '<w:r><w:t xml:space="preserve">INSERT INTO COMMENTS (USERS_ID, POSTS_ID, &quot;comment&quot;, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>comments_id</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">) VALUES (81, 4, ''no </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>esta</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> :('', 1);</w:t></w:r>'
Set-ParaXml 42 $xml

# 43: "Crear conexiones"
$xml = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Crear</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>conexiones</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
Set-ParaXml 43 $xml

# 46: "Ver trabajo especifico"
$xml = '<w:r><w:t xml:space="preserve">Ver </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>trabajo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>especifico</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
Set-ParaXml 46 $xml

# 47: "select jobs.*, users.name as username, users.photo  from jobs inner join users on users.id = jobs.users_id;"
$xml = '<w:r><w:t xml:space="preserve">select jobs.*, users.name as username, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>users.photo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">  from jobs inner join users on users.id = </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>jobs.users_id</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>;</w:t></w:r>'
Set-ParaXml 47 $xml

# New content appended at the end of the document:
#   <empty paragraph>
#   SELECT * FROM skills fetch first 10 rows only;
#   INSERT INTO skills VALUES (1, 'INGENIERIA');
Add-ParaXml ''

$xml = '<w:r><w:t>SELECT * FROM skills fetch first 10 rows only;</w:t></w:r>'
Add-ParaXml $xml

$xml = '<w:r><w:t>INSERT INTO skills VALUES (1, ''INGENIERIA'');</w:t></w:r>'
Add-ParaXml $xml

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
